$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header date (stored as text)
$ws.Range("B1").Value = "24/03/2023"

# Update activity values and hours
$ws.Range("B2").Value = 601.9
$ws.Range("C2").Value = 10

$ws.Range("B3").Value = 79
$ws.Range("C3").Value = 10

$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 10

$ws.Range("B5").Value = 558
$ws.Range("C5").Value = 10

$ws.Range("B6").Value = 366
$ws.Range("C6").Value = 10

$ws.Range("B7").Value = 115
$ws.Range("C7").Value = 10

$ws.Range("B8").Value = 151
$ws.Range("C8").Value = 10

$ws.Range("B9").Value = 392
$ws.Range("C9").Value = 10

$ws.Range("B10").Value = 47
$ws.Range("C10").Value = 10

$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 10

$ws.Range("B12").Value = 34
$ws.Range("C12").Value = 10
